$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Range("H107").Value = 601.9
$ws.Range("I107").Value = 622.82355
$ws.Range("K107").Value = 622.82355
$ws.Range("M107").Value = 1297.17645

# Row 116
$ws.Range("H116").Value = 1788.5
$ws.Range("I116").Value = 1687.5217
$ws.Range("J116").Value = 1999.6364
$ws.Range("K116").Value = 1687.5217
$ws.Range("L116").Value = 1999.6364
$ws.Range("M116").Value = 1754.4783
$ws.Range("N116").Value = -8883.636399999999

# Row 129
$ws.Range("H129").Value = 914.3958
$ws.Range("J129").Value = 953.2
$ws.Range("L129").Value = 2859.6
$ws.Range("N129").Value = -12859.6

# Row 134
$ws.Range("H134").Value = 67583.75
$ws.Range("J134").Value = 67583.75
$ws.Range("L134").Value = 67583.75
$ws.Range("N134").Value = -77723.75

# Row 135
$ws.Range("H135").Value = 100001590
$ws.Range("I135").Value = 55556190
$ws.Range("J135").Value = 166669700
$ws.Range("K135").Value = 500005710
$ws.Range("L135").Value = 1500027300
$ws.Range("M135").Value = -500003175
$ws.Range("N135").Value = -1500032370

# Row 137
$ws.Range("H137").Value = 587328.8
$ws.Range("I137").Value = 2914.9524
$ws.Range("J137").Value = 928236.9
$ws.Range("K137").Value = 8744.8572
$ws.Range("L137").Value = 2784710.7
$ws.Range("M137").Value = -6194.8572
$ws.Range("N137").Value = -2789810.7

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22100.904
$ws.Range("I32").Value = 23219.146
$ws.Range("J32").Value = 8682
$ws.Range("K32").Value = 23219.146
$ws.Range("L32").Value = 8682
$ws.Range("M32").Value = -22932.146
$ws.Range("N32").Value = -9256

$ws = $wb.Worksheets.Item("BSM")
# Row 54
$ws.Range("H54").Value = 4595.273
$ws.Range("I54").Value = 1718.5
$ws.Range("J54").Value = 12266.667
$ws.Range("K54").Value = 1718.5
$ws.Range("L54").Value = 12266.667
$ws.Range("M54").Value = -1234.5
$ws.Range("N54").Value = -13234.667

# Row 61
$ws.Range("H61").Value = 36000
$ws.Range("J61").Value = 36000
$ws.Range("L61").Value = 36000
$ws.Range("N61").Value = -36626

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 691399.9399999999
$ws.Range("I31").Value = 13144.8125
$ws.Range("J31").Value = 976981
$ws.Range("K31").Value = 13144.8125
$ws.Range("L31").Value = 976981
$ws.Range("M31").Value = -12849.8125
$ws.Range("N31").Value = -977571

# Row 34
$ws.Range("H34").Value = 691399.9399999999
$ws.Range("I34").Value = 13144.8125
$ws.Range("J34").Value = 976981
$ws.Range("K34").Value = 13144.8125
$ws.Range("L34").Value = 976981
$ws.Range("M34").Value = -12942.8125
$ws.Range("N34").Value = -977385

# Row 134
$ws.Range("H134").Value = 2802.6206
$ws.Range("I134").Value = 2249.8948
$ws.Range("J134").Value = 3852.8
$ws.Range("K134").Value = 6749.6844
$ws.Range("L134").Value = 11558.4
$ws.Range("M134").Value = -4214.6844
$ws.Range("N134").Value = -16628.4

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 12828651
$ws.Range("I5").Value = 453.55554
$ws.Range("J5").Value = 41692096
$ws.Range("K5").Value = 1360.66662
$ws.Range("L5").Value = 125076288
$ws.Range("M5").Value = -1248.66662
$ws.Range("N5").Value = -125076512

# Row 68
$ws.Range("H68").Value = 3244.2754
$ws.Range("J68").Value = 5096.9062
$ws.Range("L68").Value = 15290.7186
$ws.Range("N68").Value = -16912.7186

# Row 71
$ws.Range("H71").Value = 3244.2754
$ws.Range("J71").Value = 5096.9062
$ws.Range("L71").Value = 45872.1558
$ws.Range("N71").Value = -53984.1558

# Row 126
$ws.Range("H126").Value = 3085.5293
$ws.Range("I126").Value = 1419.8334
$ws.Range("J126").Value = 3994.0908
$ws.Range("K126").Value = 4259.5002
$ws.Range("L126").Value = 11982.2724
$ws.Range("M126").Value = 680.4997999999996
$ws.Range("N126").Value = -21862.2724

# Row 131
$ws.Range("H131").Value = 1310.1154
$ws.Range("J131").Value = 1173.8605
$ws.Range("L131").Value = 3521.5815
$ws.Range("N131").Value = -13601.5815

# Row 133
$ws.Range("H133").Value = 3124.524
$ws.Range("J133").Value = 4785.625
$ws.Range("L133").Value = 14356.875
$ws.Range("N133").Value = -24476.875

# Row 135
$ws.Range("H135").Value = 12828651
$ws.Range("I135").Value = 453.55554
$ws.Range("J135").Value = 41692096
$ws.Range("K135").Value = 4081.99986
$ws.Range("L135").Value = 375228864
$ws.Range("M135").Value = -1546.99986
$ws.Range("N135").Value = -375233934

# Row 137
$ws.Range("H137").Value = 25012.285
$ws.Range("I137").Value = 1118.8422
$ws.Range("J137").Value = 252000
$ws.Range("K137").Value = 3356.5266
$ws.Range("L137").Value = 756000
$ws.Range("M137").Value = 1743.4734
$ws.Range("N137").Value = -766200

# Row 141
$ws.Range("H141").Value = 2613.4119
$ws.Range("I141").Value = 1892.1428
$ws.Range("J141").Value = 5979.3335
$ws.Range("K141").Value = 5676.428400000001
$ws.Range("L141").Value = 17938.0005
$ws.Range("M141").Value = -496.4284000000007
$ws.Range("N141").Value = -28298.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3933
$ws.Range("I7").Value = 3933
$ws.Range("K7").Value = 3933
$ws.Range("M7").Value = -3821

# Row 46
$ws.Range("H46").Value = 1175.1666
$ws.Range("I46").Value = 762.75
$ws.Range("K46").Value = 762.75
$ws.Range("M46").Value = -574.75

# Row 55
$ws.Range("H55").Value = 131.1
$ws.Range("I55").Value = 126.833336
$ws.Range("J55").Value = 137.5
$ws.Range("K55").Value = 126.833336
$ws.Range("L55").Value = 137.5
$ws.Range("M55").Value = 46.166664
$ws.Range("N55").Value = -483.5

# Row 60
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()

# Row 68
$ws.Range("H68").Value = 4446.3076
$ws.Range("I68").Value = 3200
$ws.Range("J68").Value = 4820.2
$ws.Range("K68").Value = 3200
$ws.Range("L68").Value = 4820.2
$ws.Range("M68").Value = -2451
$ws.Range("N68").Value = -6318.2

# Row 71
$ws.Range("H71").Value = 4446.3076
$ws.Range("I71").Value = 3200
$ws.Range("J71").Value = 4820.2
$ws.Range("K71").Value = 16000
$ws.Range("L71").Value = 24101
$ws.Range("M71").Value = -12256
$ws.Range("N71").Value = -31589

# Row 126
$ws.Range("H126").Value = 3933
$ws.Range("I126").Value = 3933
$ws.Range("K126").Value = 11799
$ws.Range("M126").Value = -9329
